# Insert two new weekly price rows (Terminal Hortofrutícola Agro Chillán - Limón)
# right above the existing "1a amarillo / 2a amarillo" record that used to sit at
# rows 742-743, pushing every row from 742 downward down by two rows.
#
# Net effect matches the target diff:
#   - dimension grows from A1:T797 to A1:T799
#   - old rows 742..797 become rows 744..799 (unchanged, just shifted)
#   - two brand-new rows 742 and 743 are created, re-using the same
#     categorical columns (A,B,C,E-L,Q,R,T) as the row that used to be at
#     742/743, but carrying a new date and new price figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push everything from row 742 down by two rows.
$ws.Rows("742:743").Insert()

# 2) Seed the two freshly-inserted (blank) rows with the same categorical
#    data as the record that is now sitting at 744:745 (what used to be
#    742:743 before the insert).
$ws.Range("A744:T745").Copy()
$ws.Range("A742").PasteSpecial()

# 3) Overwrite the numeric/date fields of the new rows with their final
#    values.
$ws.Range("D742").Value = 44826
$ws.Range("M742").Value = 160
$ws.Range("N742").Value = 4500
$ws.Range("O742").Value = 5000
$ws.Range("P742").Value = 4750
$ws.Range("S742").Value = 297

$ws.Range("D743").Value = 44826
$ws.Range("M743").Value = 80
$ws.Range("N743").Value = 4000
$ws.Range("O743").Value = 4000
$ws.Range("P743").Value = 4000
$ws.Range("S743").Value = 250
